$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.018.69'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.32%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.907.99'
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.57%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.24%  '

$ws.Range("E6").Value = '  -0.50%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4831'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.29%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3807'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07364'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9338'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.82'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.37%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07792'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.899.04'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.10%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.499'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.06%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.651'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.82%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.89'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.006'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.39%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008893'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.005'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '28.053.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.35%  '

$ws.Range("E21").Value = '  +0.58%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.160'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.163.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.61%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.26'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.28%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.918'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.53'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.42%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.121'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.96%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.32'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.74%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.987'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.16%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08947'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.275'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.254'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.91%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7713'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.95%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.666'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.71%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.604'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.92%  '

$ws.Range("E37").Value = '  +0.98%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.110'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5525'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05294'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.998'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.35%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.994'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.39%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1528'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.37%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.493'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.98%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '110.43'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.55%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.71'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4829'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.83%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.006'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.50%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.648'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.61%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '68.27'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.46%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06085'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.04%  '
